# Auto-generated by analysis of the diff; sets literal (non-formula) numeric values
# for the affected currentAveragePrice / Leve profitability cells across 8 sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 1743.7028
$ws.Cells.Item(135, 9).Value = 356.6
$ws.Cells.Item(135, 10).Value = 26018
$ws.Cells.Item(135, 11).Value = 3209.4
$ws.Cells.Item(135, 12).Value = 234162
$ws.Cells.Item(135, 13).Value = -674.4000000000001
$ws.Cells.Item(135, 14).Value = -239232
$ws.Cells.Item(137, 8).Value = 13890976
$ws.Cells.Item(137, 9).Value = 21741364
$ws.Cells.Item(137, 10).Value = 1827.1538
$ws.Cells.Item(137, 11).Value = 65224092
$ws.Cells.Item(137, 12).Value = 5481.4614
$ws.Cells.Item(137, 13).Value = -65221542
$ws.Cells.Item(137, 14).Value = -10581.4614
$ws.Cells.Item(138, 8).Value = 1867.3684
$ws.Cells.Item(138, 9).Value = 956.129
$ws.Cells.Item(138, 10).Value = 2953.8462
$ws.Cells.Item(138, 11).Value = 2868.387
$ws.Cells.Item(138, 12).Value = 8861.5386
$ws.Cells.Item(138, 13).Value = 2271.613
$ws.Cells.Item(138, 14).Value = -19141.5386
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2963.923
$ws.Cells.Item(61, 9).Value = 1780.4
$ws.Cells.Item(61, 10).Value = 4577.8184
$ws.Cells.Item(61, 11).Value = 1780.4
$ws.Cells.Item(61, 12).Value = 4577.8184
$ws.Cells.Item(61, 13).Value = -1568.4
$ws.Cells.Item(61, 14).Value = -5001.8184
$ws.Cells.Item(74, 8).Value = 13160614
$ws.Cells.Item(74, 9).Value = 17244328
$ws.Cells.Item(74, 11).Value = 17244328
$ws.Cells.Item(74, 13).Value = -17243454
$ws.Cells.Item(77, 8).Value = 13160614
$ws.Cells.Item(77, 9).Value = 17244328
$ws.Cells.Item(77, 11).Value = 86221640
$ws.Cells.Item(77, 13).Value = -86217272
$ws.Cells.Item(132, 8).Value = 5528.971
$ws.Cells.Item(132, 9).Value = 2452.476
$ws.Cells.Item(132, 10).Value = 10143.714
$ws.Cells.Item(132, 11).Value = 7357.428
$ws.Cells.Item(132, 12).Value = 30431.142
$ws.Cells.Item(132, 13).Value = -4827.428
$ws.Cells.Item(132, 14).Value = -35491.142
$ws.Cells.Item(136, 8).Value = 2963.923
$ws.Cells.Item(136, 9).Value = 1780.4
$ws.Cells.Item(136, 10).Value = 4577.8184
$ws.Cells.Item(136, 11).Value = 5341.200000000001
$ws.Cells.Item(136, 12).Value = 13733.4552
$ws.Cells.Item(136, 13).Value = -2791.200000000001
$ws.Cells.Item(136, 14).Value = -18833.4552
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 1863.7273
$ws.Cells.Item(22, 9).Value = 1549.9
$ws.Cells.Item(22, 10).Value = 5002
$ws.Cells.Item(22, 11).Value = 1549.9
$ws.Cells.Item(22, 12).Value = 5002
$ws.Cells.Item(22, 13).Value = -1376.9
$ws.Cells.Item(22, 14).Value = -5348
$ws.Cells.Item(134, 8).Value = 5143.4883
$ws.Cells.Item(134, 9).Value = 2565.6667
$ws.Cells.Item(134, 10).Value = 8399.684999999999
$ws.Cells.Item(134, 11).Value = 7697.000100000001
$ws.Cells.Item(134, 12).Value = 25199.055
$ws.Cells.Item(134, 13).Value = -5162.000100000001
$ws.Cells.Item(134, 14).Value = -30269.055
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 483.66666
$ws.Cells.Item(22, 9).Value = 418.2
$ws.Cells.Item(22, 10).Value = 543.1818
$ws.Cells.Item(22, 11).Value = 418.2
$ws.Cells.Item(22, 12).Value = 543.1818
$ws.Cells.Item(22, 13).Value = -68.19999999999999
$ws.Cells.Item(22, 14).Value = -1243.1818
$ws.Cells.Item(31, 8).Value = 1746.3541
$ws.Cells.Item(31, 9).Value = 1398.3695
$ws.Cells.Item(31, 10).Value = 9750
$ws.Cells.Item(31, 11).Value = 1398.3695
$ws.Cells.Item(31, 12).Value = 9750
$ws.Cells.Item(31, 13).Value = -1103.3695
$ws.Cells.Item(31, 14).Value = -10340
$ws.Cells.Item(34, 8).Value = 1746.3541
$ws.Cells.Item(34, 9).Value = 1398.3695
$ws.Cells.Item(34, 10).Value = 9750
$ws.Cells.Item(34, 11).Value = 1398.3695
$ws.Cells.Item(34, 12).Value = 9750
$ws.Cells.Item(34, 13).Value = -1196.3695
$ws.Cells.Item(34, 14).Value = -10154
$ws.Cells.Item(58, 8).Value = 1068336.9
$ws.Cells.Item(58, 9).Value = 2234.4443
$ws.Cells.Item(58, 10).Value = 2507575
$ws.Cells.Item(58, 11).Value = 2234.4443
$ws.Cells.Item(58, 12).Value = 2507575
$ws.Cells.Item(58, 13).Value = -2031.4443
$ws.Cells.Item(58, 14).Value = -2507981
$ws.Cells.Item(107, 8).Value = 1354.56
$ws.Cells.Item(107, 9).Value = 548.94446
$ws.Cells.Item(107, 10).Value = 3426.1428
$ws.Cells.Item(107, 11).Value = 548.94446
$ws.Cells.Item(107, 12).Value = 3426.1428
$ws.Cells.Item(107, 13).Value = 1371.05554
$ws.Cells.Item(107, 14).Value = -7266.1428
$ws.Cells.Item(132, 8).Value = 2856.375
$ws.Cells.Item(132, 9).Value = 1270.4
$ws.Cells.Item(132, 10).Value = 5499.6665
$ws.Cells.Item(132, 11).Value = 3811.2
$ws.Cells.Item(132, 12).Value = 16498.9995
$ws.Cells.Item(132, 13).Value = -1281.2
$ws.Cells.Item(132, 14).Value = -21558.9995
$ws.Cells.Item(134, 8).Value = 2949.158
$ws.Cells.Item(134, 9).Value = 1369.5555
$ws.Cells.Item(134, 10).Value = 4370.8
$ws.Cells.Item(134, 11).Value = 4108.666499999999
$ws.Cells.Item(134, 12).Value = 13112.4
$ws.Cells.Item(134, 13).Value = -1573.666499999999
$ws.Cells.Item(134, 14).Value = -18182.4
$ws.Cells.Item(136, 8).Value = 1068336.9
$ws.Cells.Item(136, 9).Value = 2234.4443
$ws.Cells.Item(136, 10).Value = 2507575
$ws.Cells.Item(136, 11).Value = 6703.3329
$ws.Cells.Item(136, 12).Value = 7522725
$ws.Cells.Item(136, 13).Value = -4153.3329
$ws.Cells.Item(136, 14).Value = -7527825
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1323.05
$ws.Cells.Item(5, 10).Value = 2872.375
$ws.Cells.Item(5, 12).Value = 8617.125
$ws.Cells.Item(5, 14).Value = -8841.125
$ws.Cells.Item(122, 8).Value = 3428.7144
$ws.Cells.Item(122, 9).Value = 384.55554
$ws.Cells.Item(122, 10).Value = 3736.5505
$ws.Cells.Item(122, 11).Value = 3460.99986
$ws.Cells.Item(122, 12).Value = 33628.9545
$ws.Cells.Item(122, 13).Value = -1010.99986
$ws.Cells.Item(122, 14).Value = -38528.9545
$ws.Cells.Item(135, 8).Value = 1323.05
$ws.Cells.Item(135, 10).Value = 2872.375
$ws.Cells.Item(135, 12).Value = 25851.375
$ws.Cells.Item(135, 14).Value = -30921.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1986536.6
$ws.Cells.Item(132, 9).Value = 3207328.5
$ws.Cells.Item(132, 10).Value = 2749.75
$ws.Cells.Item(132, 11).Value = 9621985.5
$ws.Cells.Item(132, 12).Value = 8249.25
$ws.Cells.Item(132, 13).Value = -9619455.5
$ws.Cells.Item(132, 14).Value = -13309.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 27302.19
$ws.Cells.Item(55, 9).Value = 246.42105
$ws.Cells.Item(55, 11).Value = 246.42105
$ws.Cells.Item(55, 13).Value = -73.42105000000001
$ws.Cells.Item(68, 8).Value = 2576.2632
$ws.Cells.Item(68, 9).Value = 1967.375
$ws.Cells.Item(68, 11).Value = 1967.375
$ws.Cells.Item(68, 13).Value = -1218.375
$ws.Cells.Item(71, 8).Value = 2576.2632
$ws.Cells.Item(71, 9).Value = 1967.375
$ws.Cells.Item(71, 11).Value = 9836.875
$ws.Cells.Item(71, 13).Value = -6092.875
$ws.Cells.Item(100, 8).Value = 3227.15
$ws.Cells.Item(100, 9).Value = 3019
$ws.Cells.Item(100, 10).Value = 3481.5557
$ws.Cells.Item(100, 11).Value = 3019
$ws.Cells.Item(100, 12).Value = 3481.5557
$ws.Cells.Item(100, 13).Value = -2478
$ws.Cells.Item(100, 14).Value = -4563.5557
$ws.Cells.Item(132, 8).Value = 78893.28999999999
$ws.Cells.Item(132, 9).Value = 134813.5
$ws.Cells.Item(132, 10).Value = 4333
$ws.Cells.Item(132, 11).Value = 404440.5
$ws.Cells.Item(132, 12).Value = 12999
$ws.Cells.Item(132, 13).Value = -401910.5
$ws.Cells.Item(132, 14).Value = -18059
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1629.4
$ws.Cells.Item(81, 9).Value = 3467
$ws.Cells.Item(81, 10).Value = 841.8570999999999
$ws.Cells.Item(81, 11).Value = 6934
$ws.Cells.Item(81, 12).Value = 1683.7142
$ws.Cells.Item(81, 13).Value = -5873
$ws.Cells.Item(81, 14).Value = -3805.7142
$ws.Cells.Item(84, 8).Value = 1629.4
$ws.Cells.Item(84, 9).Value = 3467
$ws.Cells.Item(84, 10).Value = 841.8570999999999
$ws.Cells.Item(84, 11).Value = 34670
$ws.Cells.Item(84, 12).Value = 8418.571
$ws.Cells.Item(84, 13).Value = -29366
$ws.Cells.Item(84, 14).Value = -19026.571
$ws.Cells.Item(126, 8).Value = 4693.528
$ws.Cells.Item(126, 9).Value = 5292.067
$ws.Cells.Item(126, 10).Value = 1700.8334
$ws.Cells.Item(126, 11).Value = 15876.201
$ws.Cells.Item(126, 12).Value = 5102.5002
$ws.Cells.Item(126, 13).Value = -13406.201
$ws.Cells.Item(126, 14).Value = -10042.5002
$ws.Cells.Item(132, 8).Value = 2999.4075
$ws.Cells.Item(132, 9).Value = 2768.1538
$ws.Cells.Item(132, 10).Value = 3214.1428
$ws.Cells.Item(132, 11).Value = 8304.4614
$ws.Cells.Item(132, 12).Value = 9642.428400000001
$ws.Cells.Item(132, 13).Value = -5774.4614
$ws.Cells.Item(132, 14).Value = -14702.4284
